# Auto-generated: apply cell-value updates per Lich_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 17857268
$ws.Range("J2").Value = 35714460
$ws.Range("L2").Value = 35714460
$ws.Range("N2").Value = -35714686
$ws.Range("H17").Value = 2374662
$ws.Range("J17").Value = 2421175
$ws.Range("L17").Value = 7263525
$ws.Range("N17").Value = -7263861
$ws.Range("H19").Value = 1715.3
$ws.Range("I19").Value = 1098
$ws.Range("J19").Value = 2332.6
$ws.Range("K19").Value = 1098
$ws.Range("L19").Value = 2332.6
$ws.Range("M19").Value = -923
$ws.Range("N19").Value = -2682.6
$ws.Range("H34").Value = 2945
$ws.Range("I34").Value = 2945
$ws.Range("K34").Value = 2945
$ws.Range("M34").Value = -2742
$ws.Range("H36").Value = 2945
$ws.Range("I36").Value = 2945
$ws.Range("K36").Value = 2945
$ws.Range("M36").Value = -2230
$ws.Range("H39").Value = 307.30768
$ws.Range("I39").Value = 182.27272
$ws.Range("J39").Value = 995
$ws.Range("K39").Value = 546.81816
$ws.Range("L39").Value = 2985
$ws.Range("M39").Value = -250.81816
$ws.Range("N39").Value = -3577
$ws.Range("H41").Value = 714.4286
$ws.Range("I41").Value = 418.5
$ws.Range("J41").Value = 1109
$ws.Range("K41").Value = 418.5
$ws.Range("L41").Value = 1109
$ws.Range("M41").Value = 21.5
$ws.Range("N41").Value = -1989
$ws.Range("H64").Value = 6250.3335
$ws.Range("I64").Value = 3625.5
$ws.Range("J64").Value = 11500
$ws.Range("K64").Value = 3625.5
$ws.Range("L64").Value = 11500
$ws.Range("M64").Value = -3377.5
$ws.Range("N64").Value = -11996
$ws.Range("H67").Value = 6250.3335
$ws.Range("I67").Value = 3625.5
$ws.Range("J67").Value = 11500
$ws.Range("K67").Value = 3625.5
$ws.Range("L67").Value = 11500
$ws.Range("M67").Value = -2767.5
$ws.Range("N67").Value = -13216
$ws.Range("H76").Value = 5471.3
$ws.Range("I76").Value = 5212.8
$ws.Range("J76").Value = 5729.8
$ws.Range("K76").Value = 5212.8
$ws.Range("L76").Value = 5729.8
$ws.Range("M76").Value = -4897.8
$ws.Range("N76").Value = -6359.8
$ws.Range("H79").Value = 5471.3
$ws.Range("I79").Value = 5212.8
$ws.Range("J79").Value = 5729.8
$ws.Range("K79").Value = 5212.8
$ws.Range("L79").Value = 5729.8
$ws.Range("M79").Value = -4120.8
$ws.Range("N79").Value = -7913.8
$ws.Range("H86").Value = 3011
$ws.Range("I86").Value = 3473.3333
$ws.Range("J86").Value = 2614.7144
$ws.Range("K86").Value = 3473.3333
$ws.Range("L86").Value = 2614.7144
$ws.Range("M86").Value = -2350.3333
$ws.Range("N86").Value = -4860.7144
$ws.Range("H89").Value = 3011
$ws.Range("I89").Value = 3473.3333
$ws.Range("J89").Value = 2614.7144
$ws.Range("K89").Value = 17366.6665
$ws.Range("L89").Value = 13073.572
$ws.Range("M89").Value = -11750.6665
$ws.Range("N89").Value = -24305.572
$ws.Range("H92").Value = 3447
$ws.Range("J92").Value = 2833
$ws.Range("L92").Value = 2833
$ws.Range("N92").Value = -5329
$ws.Range("H98").Value = 2244167.8
$ws.Range("I98").Value = 2394189.8
$ws.Range("K98").Value = 2394189.8
$ws.Range("M98").Value = -2392691.8
$ws.Range("H99").Value = 333697
$ws.Range("I99").Value = 167086.17
$ws.Range("J99").Value = 666918.7
$ws.Range("K99").Value = 501258.51
$ws.Range("L99").Value = 2000756.1
$ws.Range("M99").Value = -499760.51
$ws.Range("N99").Value = -2003752.1
$ws.Range("H100").Value = 1422
$ws.Range("I100").Value = 1342.2354
$ws.Range("J100").Value = 2100
$ws.Range("K100").Value = 1342.2354
$ws.Range("L100").Value = 2100
$ws.Range("M100").Value = -801.2354
$ws.Range("N100").Value = -3182
$ws.Range("H113").Value = 7655.3335
$ws.Range("I113").Value = 9744.75
$ws.Range("K113").Value = 9744.75
$ws.Range("M113").Value = -6490.75
$ws.Range("H114").Value = 99999
$ws.Range("J114").Value = 99999
$ws.Range("L114").Value = 99999
$ws.Range("N114").Value = -108677
$ws.Range("H116").Value = 7082.5
$ws.Range("I116").Value = 7000
$ws.Range("J116").Value = 7099
$ws.Range("K116").Value = 7000
$ws.Range("L116").Value = 7099
$ws.Range("M116").Value = -3558
$ws.Range("N116").Value = -13983
$ws.Range("H122").Value = 2244167.8
$ws.Range("I122").Value = 2394189.8
$ws.Range("K122").Value = 7182569.399999999
$ws.Range("M122").Value = -7180119.399999999
$ws.Range("H129").Value = 1106.25
$ws.Range("I129").Value = 808.3333
$ws.Range("J129").Value = 2000
$ws.Range("K129").Value = 2424.9999
$ws.Range("L129").Value = 6000
$ws.Range("M129").Value = 2575.0001
$ws.Range("N129").Value = -16000
$ws.Range("H132").Value = 3674.9119
$ws.Range("I132").Value = 2431.182
$ws.Range("J132").Value = 5955.0835
$ws.Range("K132").Value = 7293.545999999999
$ws.Range("L132").Value = 17865.2505
$ws.Range("M132").Value = -4763.545999999999
$ws.Range("N132").Value = -22925.2505
$ws.Range("H135").Value = 1712.8572
$ws.Range("I135").Value = 1536
$ws.Range("J135").Value = 2155
$ws.Range("K135").Value = 13824
$ws.Range("L135").Value = 19395
$ws.Range("M135").Value = -11289
$ws.Range("N135").Value = -24465
$ws.Range("H137").Value = 23563.127
$ws.Range("I137").Value = 29318.05
$ws.Range("J137").Value = 8216.666999999999
$ws.Range("K137").Value = 87954.14999999999
$ws.Range("L137").Value = 24650.001
$ws.Range("M137").Value = -85404.14999999999
$ws.Range("N137").Value = -29750.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3814.0156
$ws.Range("I32").Value = 2889.2856
$ws.Range("K32").Value = 2889.2856
$ws.Range("M32").Value = -2602.2856
$ws.Range("H61").Value = 3003.6274
$ws.Range("I61").Value = 1899.8823
$ws.Range("K61").Value = 1899.8823
$ws.Range("M61").Value = -1687.8823
$ws.Range("H74").Value = 52283.258
$ws.Range("I74").Value = 63073.22
$ws.Range("J74").Value = 2957.7144
$ws.Range("K74").Value = 63073.22
$ws.Range("L74").Value = 2957.7144
$ws.Range("M74").Value = -62199.22
$ws.Range("N74").Value = -4705.7144
$ws.Range("H77").Value = 52283.258
$ws.Range("I77").Value = 63073.22
$ws.Range("J77").Value = 2957.7144
$ws.Range("K77").Value = 315366.1
$ws.Range("L77").Value = 14788.572
$ws.Range("M77").Value = -310998.1
$ws.Range("N77").Value = -23524.572
$ws.Range("H97").Value = 1016.05554
$ws.Range("I97").Value = 775.3077
$ws.Range("K97").Value = 775.3077
$ws.Range("M97").Value = -279.3077
$ws.Range("H110").Value = 3381.4
$ws.Range("I110").Value = 1726.75
$ws.Range("K110").Value = 1726.75
$ws.Range("M110").Value = 318.25
$ws.Range("H122").Value = 2768.1052
$ws.Range("I122").Value = 2200.4614
$ws.Range("K122").Value = 6601.3842
$ws.Range("M122").Value = -4151.3842
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H132").Value = 1736.5897
$ws.Range("I132").Value = 1697.931
$ws.Range("K132").Value = 5093.793
$ws.Range("M132").Value = -2563.793
$ws.Range("H136").Value = 3003.6274
$ws.Range("I136").Value = 1899.8823
$ws.Range("K136").Value = 5699.6469
$ws.Range("M136").Value = -3149.6469
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1700.862
$ws.Range("I20").Value = 1601.3684
$ws.Range("K20").Value = 1601.3684
$ws.Range("M20").Value = -1354.3684
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("H86").Value = 1529.1111
$ws.Range("I86").Value = 1095.875
$ws.Range("J86").Value = 4995
$ws.Range("K86").Value = 1095.875
$ws.Range("L86").Value = 4995
$ws.Range("M86").Value = 27.125
$ws.Range("N86").Value = -7241
$ws.Range("H89").Value = 1529.1111
$ws.Range("I89").Value = 1095.875
$ws.Range("J89").Value = 4995
$ws.Range("K89").Value = 5479.375
$ws.Range("L89").Value = 24975
$ws.Range("M89").Value = 136.625
$ws.Range("N89").Value = -36207
$ws.Range("H99").Value = 4452.1304
$ws.Range("J99").Value = 4129.8335
$ws.Range("L99").Value = 4129.8335
$ws.Range("N99").Value = -7125.8335
$ws.Range("H105").Value = 1948.8235
$ws.Range("I105").Value = 1953.5483
$ws.Range("J105").Value = 1900
$ws.Range("K105").Value = 1953.5483
$ws.Range("L105").Value = 1900
$ws.Range("M105").Value = -206.5482999999999
$ws.Range("N105").Value = -5394
$ws.Range("H112").Value = 19653
$ws.Range("J112").Value = 19653
$ws.Range("L112").Value = 19653
$ws.Range("N112").Value = -22607
$ws.Range("H134").Value = 3284.75
$ws.Range("I134").Value = 2393.2778
$ws.Range("J134").Value = 4889.4
$ws.Range("K134").Value = 7179.8334
$ws.Range("L134").Value = 14668.2
$ws.Range("M134").Value = -4644.8334
$ws.Range("N134").Value = -19738.2
$ws.Range("M75").ClearContents()
$ws.Range("M78").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 483.6087
$ws.Range("J7").Value = 529.7778
$ws.Range("L7").Value = 529.7778
$ws.Range("N7").Value = -755.7778
$ws.Range("H16").Value = 1434.7576
$ws.Range("I16").Value = 1470.7778
$ws.Range("J16").Value = 1391.5333
$ws.Range("K16").Value = 1470.7778
$ws.Range("L16").Value = 1391.5333
$ws.Range("M16").Value = -1183.7778
$ws.Range("N16").Value = -1965.5333
$ws.Range("H22").Value = 364.55173
$ws.Range("I22").Value = 365.7037
$ws.Range("K22").Value = 365.7037
$ws.Range("M22").Value = -15.70370000000003
$ws.Range("H31").Value = 324641.03
$ws.Range("I31").Value = 556554.75
$ws.Range("K31").Value = 556554.75
$ws.Range("M31").Value = -556259.75
$ws.Range("H34").Value = 324641.03
$ws.Range("I34").Value = 556554.75
$ws.Range("K34").Value = 556554.75
$ws.Range("M34").Value = -556352.75
$ws.Range("H58").Value = 2598.2666
$ws.Range("I58").Value = 2262.4348
$ws.Range("J58").Value = 3701.7144
$ws.Range("K58").Value = 2262.4348
$ws.Range("L58").Value = 3701.7144
$ws.Range("M58").Value = -2059.4348
$ws.Range("N58").Value = -4107.7144
$ws.Range("H62").Value = 3351.8572
$ws.Range("I62").Value = 3362.9
$ws.Range("J62").Value = 3324.25
$ws.Range("K62").Value = 3362.9
$ws.Range("L62").Value = 3324.25
$ws.Range("M62").Value = -2738.9
$ws.Range("N62").Value = -4572.25
$ws.Range("H65").Value = 3351.8572
$ws.Range("I65").Value = 3362.9
$ws.Range("J65").Value = 3324.25
$ws.Range("K65").Value = 16814.5
$ws.Range("L65").Value = 16621.25
$ws.Range("M65").Value = -13694.5
$ws.Range("N65").Value = -22861.25
$ws.Range("H74").Value = 64307
$ws.Range("J74").Value = 67812.5
$ws.Range("L74").Value = 67812.5
$ws.Range("N74").Value = -69560.5
$ws.Range("H77").Value = 64307
$ws.Range("J77").Value = 67812.5
$ws.Range("L77").Value = 203437.5
$ws.Range("N77").Value = -212173.5
$ws.Range("H99").Value = 302520.12
$ws.Range("I99").Value = 779275.9399999999
$ws.Range("K99").Value = 779275.9399999999
$ws.Range("M99").Value = -777777.9399999999
$ws.Range("H113").Value = 1434.7576
$ws.Range("I113").Value = 1470.7778
$ws.Range("J113").Value = 1391.5333
$ws.Range("K113").Value = 1470.7778
$ws.Range("L113").Value = 1391.5333
$ws.Range("M113").Value = 699.2221999999999
$ws.Range("N113").Value = -5731.5333
$ws.Range("H126").Value = 302520.12
$ws.Range("I126").Value = 779275.9399999999
$ws.Range("K126").Value = 2337827.82
$ws.Range("M126").Value = -2335357.82
$ws.Range("H134").Value = 6496.613
$ws.Range("I134").Value = 6774.074
$ws.Range("K134").Value = 20322.222
$ws.Range("M134").Value = -17787.222
$ws.Range("H136").Value = 2598.2666
$ws.Range("I136").Value = 2262.4348
$ws.Range("J136").Value = 3701.7144
$ws.Range("K136").Value = 6787.3044
$ws.Range("L136").Value = 11105.1432
$ws.Range("M136").Value = -4237.3044
$ws.Range("N136").Value = -16205.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 5000
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("H17").Value = 8900
$ws.Range("J17").Value = 8900
$ws.Range("L17").Value = 26700
$ws.Range("N17").Value = -27038
$ws.Range("H37").Value = 100112870
$ws.Range("J37").Value = 100112870
$ws.Range("L37").Value = 300338610
$ws.Range("N37").Value = -300338834
$ws.Range("H39").Value = 1153.5454
$ws.Range("I39").Value = 336.625
$ws.Range("J39").Value = 3332
$ws.Range("K39").Value = 1009.875
$ws.Range("L39").Value = 9996
$ws.Range("M39").Value = -715.875
$ws.Range("N39").Value = -10584
$ws.Range("H41").Value = 900.3
$ws.Range("J41").Value = 1433.6666
$ws.Range("L41").Value = 4300.9998
$ws.Range("N41").Value = -4976.9998
$ws.Range("H44").Value = 2034.4
$ws.Range("I44").Value = 720.6667
$ws.Range("J44").Value = 2773.375
$ws.Range("K44").Value = 2162.0001
$ws.Range("L44").Value = 8320.125
$ws.Range("M44").Value = -1764.0001
$ws.Range("N44").Value = -9116.125
$ws.Range("H98").Value = 397.5
$ws.Range("I98").Value = 441.25
$ws.Range("J98").Value = 368.33334
$ws.Range("K98").Value = 1323.75
$ws.Range("L98").Value = 1105.00002
$ws.Range("M98").Value = 174.25
$ws.Range("N98").Value = -4101.000019999999
$ws.Range("H131").Value = 1440.9487
$ws.Range("J131").Value = 1452.5526
$ws.Range("L131").Value = 4357.6578
$ws.Range("N131").Value = -14437.6578
$ws.Range("H140").Value = 24589.555
$ws.Range("I140").Value = 34107
$ws.Range("J140").Value = 5554.6665
$ws.Range("K140").Value = 102321
$ws.Range("L140").Value = 16663.9995
$ws.Range("M140").Value = -97141
$ws.Range("N140").Value = -27023.9995
$ws.Range("M9").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 51
$ws.Range("J2").Value = 37
$ws.Range("L2").Value = 37
$ws.Range("N2").Value = -263
$ws.Range("H97").Value = 2935.4666
$ws.Range("J97").Value = 4735.25
$ws.Range("L97").Value = 4735.25
$ws.Range("N97").Value = -5727.25
$ws.Range("H102").Value = 3103.875
$ws.Range("I102").Value = 3065.7144
$ws.Range("J102").Value = 3371
$ws.Range("K102").Value = 3065.7144
$ws.Range("L102").Value = 3371
$ws.Range("M102").Value = -1443.7144
$ws.Range("N102").Value = -6615
$ws.Range("H113").Value = 11666.267
$ws.Range("I113").Value = 12120.091
$ws.Range("J113").Value = 10418.25
$ws.Range("K113").Value = 12120.091
$ws.Range("L113").Value = 10418.25
$ws.Range("M113").Value = -9950.091
$ws.Range("N113").Value = -14758.25
$ws.Range("H123").Value = 39460.8
$ws.Range("J123").Value = 39460.8
$ws.Range("L123").Value = 39460.8
$ws.Range("N123").Value = -44360.8
$ws.Range("H132").Value = 28310.512
$ws.Range("I132").Value = 33936.727
$ws.Range("J132").Value = 5102.375
$ws.Range("K132").Value = 101810.181
$ws.Range("L132").Value = 15307.125
$ws.Range("M132").Value = -99280.181
$ws.Range("N132").Value = -20367.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3118.8333
$ws.Range("J16").Value = 27994
$ws.Range("L16").Value = 27994
$ws.Range("N16").Value = -28334
$ws.Range("H22").Value = 1112.3636
$ws.Range("I22").Value = 888.2
$ws.Range("J22").Value = 1209.826
$ws.Range("K22").Value = 888.2
$ws.Range("L22").Value = 1209.826
$ws.Range("M22").Value = -593.2
$ws.Range("N22").Value = -1799.826
$ws.Range("H27").Value = 1112.3636
$ws.Range("I27").Value = 888.2
$ws.Range("J27").Value = 1209.826
$ws.Range("K27").Value = 888.2
$ws.Range("L27").Value = 1209.826
$ws.Range("M27").Value = -781.2
$ws.Range("N27").Value = -1423.826
$ws.Range("H38").Value = 73748
$ws.Range("J38").Value = 99999
$ws.Range("L38").Value = 99999
$ws.Range("N38").Value = -100819
$ws.Range("H40").Value = 24527
$ws.Range("I40").Value = 26309.8
$ws.Range("J40").Value = 6699
$ws.Range("K40").Value = 26309.8
$ws.Range("L40").Value = 6699
$ws.Range("M40").Value = -26173.8
$ws.Range("N40").Value = -6971
$ws.Range("H41").Value = 49999
$ws.Range("I41").Value = 49999
$ws.Range("K41").Value = 49999
$ws.Range("M41").Value = -49561
$ws.Range("H55").Value = 4750.1353
$ws.Range("I55").Value = 724.12
$ws.Range("J55").Value = 13137.667
$ws.Range("K55").Value = 724.12
$ws.Range("L55").Value = 13137.667
$ws.Range("M55").Value = -551.12
$ws.Range("N55").Value = -13483.667
$ws.Range("H61").Value = 1253.36
$ws.Range("I61").Value = 1177
$ws.Range("K61").Value = 1177
$ws.Range("M61").Value = -975
$ws.Range("H68").Value = 3852.2856
$ws.Range("J68").Value = 2400
$ws.Range("L68").Value = 2400
$ws.Range("N68").Value = -3898
$ws.Range("H71").Value = 3852.2856
$ws.Range("J71").Value = 2400
$ws.Range("L71").Value = 12000
$ws.Range("N71").Value = -19488
$ws.Range("H81").Value = 37500
$ws.Range("J81").Value = 37500
$ws.Range("L81").Value = 37500
$ws.Range("N81").Value = -39496
$ws.Range("H84").Value = 37500
$ws.Range("J84").Value = 37500
$ws.Range("L84").Value = 112500
$ws.Range("N84").Value = -122484
$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990
$ws.Range("H113").Value = 1253.36
$ws.Range("I113").Value = 1177
$ws.Range("K113").Value = 1177
$ws.Range("M113").Value = 993
$ws.Range("H122").Value = 5319.4
$ws.Range("I122").Value = 6900
$ws.Range("J122").Value = 4924.25
$ws.Range("K122").Value = 20700
$ws.Range("L122").Value = 14772.75
$ws.Range("M122").Value = -18250
$ws.Range("N122").Value = -19672.75
$ws.Range("H136").Value = 1229.0869
$ws.Range("I136").Value = 1039.05
$ws.Range("K136").Value = 3117.15
$ws.Range("M136").Value = -567.1499999999996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 4499
$ws.Range("I26").Value = 4499
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 4499
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -4206
$ws.Range("H31").Value = 18499.5
$ws.Range("J31").Value = 18499.5
$ws.Range("L31").Value = 18499.5
$ws.Range("N31").Value = -19195.5
$ws.Range("H51").Value = 20000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("H107").Value = 709.4
$ws.Range("I107").Value = 517.8333
$ws.Range("J107").Value = 996.75
$ws.Range("K107").Value = 1553.4999
$ws.Range("L107").Value = 2990.25
$ws.Range("M107").Value = 366.5001
$ws.Range("N107").Value = -6830.25
$ws.Range("H122").Value = 25834.75
$ws.Range("I122").Value = 2573.9
$ws.Range("J122").Value = 142139
$ws.Range("K122").Value = 7721.700000000001
$ws.Range("L122").Value = 426417
$ws.Range("M122").Value = -5271.700000000001
$ws.Range("N122").Value = -431317
$ws.Range("H132").Value = 2298.4
$ws.Range("I132").Value = 2249.5
$ws.Range("K132").Value = 6748.5
$ws.Range("M132").Value = -4218.5
$ws.Range("H136").Value = 371285.53
$ws.Range("I136").Value = 385531.97
$ws.Range("J136").Value = 878
$ws.Range("K136").Value = 1156595.91
$ws.Range("L136").Value = 2634
$ws.Range("M136").Value = -1154045.91
$ws.Range("N136").Value = -7734
$ws.Range("N26").ClearContents()
$ws.Range("N51").ClearContents()
